# Swap the data (columns B through AB) between paired match rows.
# The "id" column (A) stays fixed per row; everything else (B:AB) for each
# pair of rows gets exchanged, effectively swapping which physical row each
# match's data lives on while keeping the row's sequence id unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowPairs = @(
    @(13, 14),
    @(20, 21),
    @(26, 27),
    @(28, 29),
    @(38, 39),
    @(43, 44),
    @(47, 48),
    @(54, 55),
    @(56, 57)
)

# Columns B (2) through AB (28) inclusive
$firstCol = 2
$lastCol = 28

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    for ($col = $firstCol; $col -le $lastCol; $col++) {
        $cell1 = $ws.Cells.Item($r1, $col)
        $cell2 = $ws.Cells.Item($r2, $col)

        $v1 = $cell1.Value2
        $v2 = $cell2.Value2

        $cell1.Value2 = $v2
        $cell2.Value2 = $v1
    }
}
